$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename DB entity/column names from PascalCase/class-style to snake_case,
# matching the latest DDL files (table + column names).
$ws.Range("B1").Value = 'sdc_observation'
$ws.Range("C1").Value = 'template_instance (DiagReport)'
$ws.Range("D1").Value = 'template_sdc'
$ws.Range("E1").Value = 'template_term_map'
$ws.Range("F1").Value = 'template_map_content'
$ws.Range("B2").Value = 'sdc_observation_id'
$ws.Range("C2").Value = 'template_instance_id'
$ws.Range("D2").Value = 'template_sdc_id'
$ws.Range("E2").Value = 'template_term_map_id'
$ws.Range("F2").Value = 'template_map_content_id'
$ws.Range("B3").Value = 'template_instance_id'
$ws.Range("C3").Value = 'template_instance_version_guid'
$ws.Range("D3").Value = 'sdc_form_design_sdcid (FD or Pkg)'
$ws.Range("E3").Value = 'template_map_sdcid'
$ws.Range("F3").Value = 'template_term_map_id'
$ws.Range("B4").Value = 'parent_sdc_observation_id'
$ws.Range("C4").Value = 'template_instance_version_uri'
$ws.Range("D4").Value = 'base_uri'
$ws.Range("E4").Value = 'template'
$ws.Range("F4").Value = 'target_id'
$ws.Range("B5").Value = 'parent_instance_guid'
$ws.Range("C5").Value = 'template_sdc_id'
$ws.Range("D5").Value = 'lineage'
$ws.Range("E5").Value = 'template_sdc_id'
$ws.Range("F5").Value = 'code'
$ws.Range("C6").Value = 'instance_version_date'
$ws.Range("D6").Value = 'version'
$ws.Range("E6").Value = 'map_xml'
$ws.Range("F6").Value = 'code_text'
$ws.Range("B7").Value = 'section_sdcid'
$ws.Range("C7").Value = 'diag_report_props…'
$ws.Range("D7").Value = 'full_uri'
$ws.Range("E7").Value = 'code_system_name'
$ws.Range("F7").Value = 'code_match'
$ws.Range("B8").Value = 'section_guid'
$ws.Range("D8").Value = 'form_title'
$ws.Range("E8").Value = 'code_system_release_date'
$ws.Range("C9").Value = 'surg_path_sdcid'
$ws.Range("D9").Value = 'sdc_xml'
$ws.Range("E9").Value = 'code_system_version'
$ws.Range("B10").Value = 'question_text'
$ws.Range("D10").Value = 'doc_type (FD or Pkg)'
$ws.Range("E10").Value = 'code_system_oid'
$ws.Range("B11").Value = 'question_instance_guid'
$ws.Range("C11").Value = 'person_id'
$ws.Range("E11").Value = 'code_system_uri'
$ws.Range("B12").Value = 'question_sdcid'
$ws.Range("C12").Value = 'visit_occurrence_id'
$ws.Range("B13").Value = 'list_item_text'
$ws.Range("C13").Value = 'provider_id'
$ws.Range("B14").Value = 'list_item_sdcid'
$ws.Range("B15").Value = 'list_item_instance_guid'
$ws.Range("C15").Value = 'report_text'
$ws.Range("B16").Value = 'list_item_parent_guid ?'
$ws.Range("B17").Value = 'response'
$ws.Range("D18").Value = 'observation_specimens'
$ws.Range("E18").Value = 'sdc_specimen'
$ws.Range("B19").Value = 'units'
$ws.Range("D19").Value = 'observation_specimens_id'
$ws.Range("E19").Value = 'sdc_specimen_id'
$ws.Range("B20").Value = 'units_system'
$ws.Range("D20").Value = 'sdc_observation_id'
$ws.Range("E20").Value = 'parent_sdc_specimen_id'
$ws.Range("B21").Value = 'data_type'
$ws.Range("D21").Value = 'sdc_specimen_id'
$ws.Range("E21").Value = 'patient_id'
$ws.Range("B22").Value = 'response_int'
$ws.Range("E22").Value = 'visit_occurrence_id'
$ws.Range("B23").Value = 'response_float'
$ws.Range("B24").Value = 'response_datetime'
$ws.Range("E24").Value = 'specimen_type_text'
$ws.Range("B25").Value = 'response_string'
$ws.Range("E25").Value = 'specimen_type_code'
$ws.Range("B26").Value = 'obs_datetime'
$ws.Range("E26").Value = 'source_site_text'
$ws.Range("B27").Value = 'sdc_order'
$ws.Range("E27").Value = 'source_site_code'
$ws.Range("B28").Value = 'sdc_repeat_level'
$ws.Range("E28").Value = 'collection_method_text'
$ws.Range("B29").Value = 'sdc_comments'
$ws.Range("E29").Value = 'collection_method_code'
$ws.Range("E30").Value = 'speciment_count'
$ws.Range("E31").Value = 'collection_date'
$ws.Range("B33").Value = 'person_id'
$ws.Range("B34").Value = 'visit_occurrence_id'
$ws.Range("B35").Value = 'provider_id'

# Update the saved view: zoom level and active selection.
$ws.Range("F21").Select()
$excel.ActiveWindow.Zoom = 185
